# Update Tac1-Tacr1.xlsx with new TPM data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 values (Sending cluster FAPs -> Target cluster ECs)
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Tac1"
$ws.Range("C2").Value = "Tacr1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.480258666666666
$ws.Range("H2").Value = 7.440776
$ws.Range("I2").Value = 0.9476581720434079
$ws.Range("J2").Value = 0.947658172043408
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.1826916666666667
$ws.Range("N2").Value = 0.548075
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.4531225895777777
$ws.Range("R2").Value = 4.0781033062
$ws.Range("S2").Value = 0.9476581720434079
$ws.Range("T2").Value = 0.947658172043408

# Row 3 values (Sending cluster Resolving-Mac -> Target cluster ECs)
$ws.Range("A3").Value = "Resolving-Mac"
$ws.Range("B3").Value = "Tac1"
$ws.Range("C3").Value = "Tacr1"
$ws.Range("D3").Value = "ECs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.1369916666666667
$ws.Range("H3").Value = 0.410975
$ws.Range("I3").Value = 0.0523418279565921
$ws.Range("J3").Value = 0.0523418279565921
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.1826916666666667
$ws.Range("N3").Value = 0.548075
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.02502723590277777
$ws.Range("R3").Value = 0.225245123125
$ws.Range("S3").Value = 0.0523418279565921
$ws.Range("T3").Value = 0.0523418279565921
